# Add the "Volumetry Efficienccy Stock" worksheet with the VE data table,
# placing it after the existing sheets and making it the active/selected tab.

$wb = $excel.ActiveWorkbook

# Reference sheet used as the source of the existing cell style (vertical
# center + wrap text) so we re-use style index 1 instead of creating a new one.
$styleSource = $wb.Worksheets.Item(1).Range("A1")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Volumetry Efficienccy Stock"

$rowsData = @(
  @(1.085,1.085,1.0609999999999999,1.0760000000000001,1.0429999999999999,1.0329999999999999,1.0549999999999999,1.056,1.004,0.97740000000000005,0.96519999999999995,1.0229999999999999,1.1100000000000001,0.99750000000000005,1.0029999999999999,1,1.006,1.02),
  @(1.085,1.085,1.0649999999999999,1.0249999999999999,1.02,1.034,1.042,1.0269999999999999,0.96799999999999997,0.97519999999999996,0.96779999999999999,1.036,1.1100000000000001,0.99750000000000005,1.0029999999999999,1,1.006,1.02),
  @(1.0669999999999999,1.0669999999999999,1.0609999999999999,1.0169999999999999,1.014,1.0209999999999999,1.042,1.018,0.97050000000000003,0.97289999999999999,0.97040000000000004,1.05,1.056,0.99750000000000005,1.0029999999999999,1,1.006,1.02),
  @(1.0660000000000001,1.0660000000000001,1.034,1.0069999999999999,1.0049999999999999,1.018,1.016,1.0009999999999999,0.98029999999999995,0.95909999999999995,0.96589999999999998,1.0269999999999999,1.0009999999999999,1.0329999999999999,1.012,1,0.995,1),
  @(1.0609999999999999,1.0609999999999999,1.022,0.99960000000000004,0.98160000000000003,1.008,0.98440000000000005,0.98839999999999995,0.97629999999999995,0.94540000000000002,0.97350000000000003,0.9879,1.01,1.0529999999999999,1.0149999999999999,0.995,1,0.98499999999999999),
  @(1.0589999999999999,1.0589999999999999,1.002,0.97850000000000004,0.97370000000000001,0.98850000000000005,0.96909999999999996,0.96630000000000005,0.97240000000000004,0.95240000000000002,0.94989999999999997,1.0169999999999999,1.0289999999999999,1.0549999999999999,1,0.99199999999999999,1,1.0249999999999999),
  @(1.0580000000000001,1.0580000000000001,0.99680000000000002,0.97119999999999995,0.97140000000000004,0.97460000000000002,0.96530000000000005,0.94430000000000003,0.95879999999999999,0.95050000000000001,0.93610000000000004,1.02,1.0149999999999999,1.03,0.98399999999999999,1.01,0.99299999999999999,1.01),
  @(1.034,1.034,0.98240000000000005,0.97599999999999998,0.97470000000000001,0.96540000000000004,0.95950000000000002,0.95550000000000002,0.94710000000000005,0.93569999999999998,0.98099999999999998,1.016,1.0049999999999999,0.99750000000000005,1.01,1,0.99,1),
  @(1.01,1.01,0.96799999999999997,0.96130000000000004,0.98360000000000003,0.96550000000000002,0.95379999999999998,0.95220000000000005,0.93989999999999996,0.98080000000000001,0.98499999999999999,1.004,0.995,1.038,1,1.004,0.98799999999999999,0.98),
  @(1,1,0.97289999999999999,0.94669999999999999,0.96679999999999999,0.94599999999999995,0.95660000000000001,0.94779999999999998,0.9546,0.96330000000000005,0.99370000000000003,0.99250000000000005,0.995,1.0229999999999999,0.99199999999999999,0.98299999999999998,0.97599999999999998,0.97),
  @(0.99,0.99,0.95499999999999996,0.94610000000000005,0.95,0.92659999999999998,0.93559999999999999,0.93930000000000002,0.95550000000000002,0.97409999999999997,1.0029999999999999,0.97,0.99750000000000005,0.99750000000000005,0.98499999999999999,0.995,0.98,0.97),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.92500000000000004,0.91610000000000003,0.92879999999999996,0.93830000000000002,0.98499999999999999,0.98750000000000004,0.99,0.98499999999999999,0.98250000000000004,0.97,0.98799999999999999,0.97,0.95),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.88849999999999996,0.91300000000000003,0.93230000000000002,1.0029999999999999,1,0.98250000000000004,0.97250000000000003,0.97,0.96,0.97599999999999998,0.95499999999999996,0.96199999999999997),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.90859999999999996,0.9375,0.98,0.98499999999999999,0.97499999999999998,0.97,0.95250000000000001,0.97250000000000003,0.96,1,0.97699999999999998,0.98),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.90859999999999996,0.9375,0.98,0.96419999999999995,0.97499999999999998,0.96630000000000005,0.94520000000000004,0.96879999999999999,0.9768,1,0.97699999999999998,0.89639999999999997),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.90859999999999996,0.9375,0.98,0.96319999999999995,0.95079999999999998,0.96499999999999997,0.9425,0.95799999999999996,0.98299999999999998,0.96850000000000003,0.93820000000000003,0.89429999999999998),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.90859999999999996,0.9375,0.98,0.96209999999999996,0.92500000000000004,0.95040000000000002,0.92789999999999995,0.94650000000000001,0.98299999999999998,0.93500000000000005,0.89700000000000002,0.89219999999999999),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.90859999999999996,0.9375,0.98,0.95950000000000002,0.91839999999999999,0.91539999999999999,0.89290000000000003,0.91900000000000004,0.93300000000000005,0.9284,0.89039999999999997,0.88690000000000002),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.90859999999999996,0.9375,0.98,0.95420000000000005,0.90500000000000003,0.90339999999999998,0.88090000000000002,0.90700000000000003,0.92100000000000004,0.91500000000000004,0.877,0.87639999999999996),
  @(0.99,0.99,0.95499999999999996,0.93500000000000005,0.91,0.90500000000000003,0.90859999999999996,0.9375,0.98,0.95420000000000005,0.90500000000000003,0.89539999999999997,0.87290000000000001,0.89900000000000002,0.91300000000000003,0.91500000000000004,0.877,0.87639999999999996)
)

$nRows = $rowsData.Length
$nCols = $rowsData[0].Length

$data = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
    $rowVals = $rowsData[$r]
    for ($c = 0; $c -lt $nCols; $c++) {
        $data[$r,$c] = $rowVals[$c]
    }
}

$targetRange = $newSheet.Range("A1").Resize($nRows, $nCols)
$targetRange.Value = $data

# Apply the same style (vertical-centered, wrapped text) used throughout the
# rest of the workbook by copying formatting from an existing styled cell.
$styleSource.Copy()
$targetRange.PasteSpecial(-4122)

# Make the new sheet the active tab and set its selection, matching the
# workbook state saved by the author.
$newSheet.Activate()
$newSheet.Range("C23").Select()
